$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 new values (swapped with row 6's original values)
$ws.Range("D2").Value = 44200
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1450
$ws.Range("P2").Value = 1450

# Row 6 new values (swapped with row 2's original values)
$ws.Range("D6").Value = 44638
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 2800
$ws.Range("M6").Value = 2650
$ws.Range("P6").Value = 2650
